$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header labels (columns D, E, F added; B/C header text replaced)
$ws.Range("B1").Value = "Contiguous not selected"
$ws.Range("C1").Value = "Contiguous selected"
$ws.Range("D1").Value = "Ineligible"
$ws.Range("E1").Value = "LIC not selected"
$ws.Range("F1").Value = "LIC selected"

# Time period labels used repeatedly
$periods = @("2014-09","2015-09","2016-09","2017-09","2018-09","2019-09","2020-09","2021-09","2022-09","2023-09","2024-09")

# Data rows (each row: B, C, D, E, F)
$data = @(
    @(15803506, 227030, 55062798, 33946627, 10670061),
    @(15960120, 230742, 55712041, 34134976, 10741003),
    @(16110051, 233418, 56370128, 34322254, 10808051),
    @(16259861, 236255, 57078682, 34515198, 10876607),
    @(16437369, 238980, 57816819, 34719942, 10956502),
    @(16616200, 241444, 58582032, 34948131, 11040777),
    @(16792486, 244689, 59313397, 35146197, 11120071),
    @(16994655, 248863, 60178302, 35393243, 11229007),
    @(17216365, 253408, 61059978, 35687306, 11368159),
    @(17394800, 256691, 61765087, 35919647, 11485528),
    @(17680903, 263739, 62780406, 36304733, 11692274)
)

# Rows 2-12: first block of time periods
for ($i = 0; $i -lt 11; $i++) {
    $r = 2 + $i
    $ws.Cells.Item($r, 1).Value = $periods[$i]
    $ws.Cells.Item($r, 2).Value = $data[$i][0]
    $ws.Cells.Item($r, 3).Value = $data[$i][1]
    $ws.Cells.Item($r, 4).Value = $data[$i][2]
    $ws.Cells.Item($r, 5).Value = $data[$i][3]
    $ws.Cells.Item($r, 6).Value = $data[$i][4]
}

# Rows 13-23: second (duplicate) block of time periods
for ($i = 0; $i -lt 11; $i++) {
    $r = 13 + $i
    $ws.Cells.Item($r, 1).Value = $periods[$i]
    $ws.Cells.Item($r, 2).Value = $data[$i][0]
    $ws.Cells.Item($r, 3).Value = $data[$i][1]
    $ws.Cells.Item($r, 4).Value = $data[$i][2]
    $ws.Cells.Item($r, 5).Value = $data[$i][3]
    $ws.Cells.Item($r, 6).Value = $data[$i][4]
}

# Row 24: Pre-OZs Avg. Annual Growth
$ws.Cells.Item(24, 1).Value = "Pre-OZs Avg. Annual Growth"
$ws.Cells.Item(24, 2).Value = 1.00799658051637
$ws.Cells.Item(24, 3).Value = 1.2389285086471
$ws.Cells.Item(24, 4).Value = 1.24679939773407
$ws.Cells.Item(24, 5).Value = 0.583211899326659
$ws.Cells.Item(24, 6).Value = 0.685427091015054

# Row 25: Post-OZs Avg. Annual Growth
$ws.Cells.Item(25, 1).Value = "Post-OZs Avg. Annual Growth"
$ws.Cells.Item(25, 2).Value = 1.25012576919903
$ws.Cells.Item(25, 3).Value = 1.7834789518827
$ws.Cells.Item(25, 4).Value = 1.3940735455782
$ws.Cells.Item(25, 5).Value = 0.764724052376205
$ws.Cells.Item(25, 6).Value = 1.15390764023835
